# ICTP/application.docx edit script
# Applies the "se termina la aplicacion a SAIFR" revision.

$d = $word.ActiveDocument

# Use wildcard=False, MatchCase=True for the Find/Execute calls.
# Execute signature (order used throughout):
# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

# ---------------------------------------------------------------------------
# 1) Collapse the "Reason for Participation (...)" heading into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Reason for Participation (maximum 4000 characters):", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Reason for Participation (maximum 4000 characters):", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Rewrite the tail of the "I aim to better understand..." paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "I also believe this is a great chance of getting immersed at the international scale with the theoretical physics community. In particular, I am excited to get to know colleagues with similar interests who may eventually turn into collaborators.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In particular, I am very excited about the lectures on topology for physicists. This is an immensely important subject in modern theoretical endeavors which is usually forgotten at the undergraduate level. On the other hand, the other set of lectures will provide me with the chance to learn more about many subjects which are in constant contact with my current research. Indeed, algebraic formulations provide a more structural understanding of topics such as field theory and phase transitions. Examples can be found in the Haag-Kastler axioms and the application of index theory to the topic of topological phases.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert a brand-new paragraph (+ trailing blank separator) right before
#    "I am also very interested in the SAIFR-Perimeter Fellowship...".
# ---------------------------------------------------------------------------
$newPara = "In the additional information and summary sections I list why I believe I have the right prerequisites to assist to this school. However, I also believe that this is a great chance of getting immersed at the international scale with the theoretical physics community. In particular, I am excited to get to know colleagues with similar interests who may eventually turn into collaborators. This is an opportunity to consolidate a network of Latin American physicists that will make our continent a force to be reckoned with at the frontiers of theoretical physics."

$saifrRange = $d.Content.Find.Execute(
    "I am also very interested in the SAIFR-Perimeter Fellowship", $true,
    $false, $false, $false, $false, $true, 1, $false, "", 0)
$saifrPara = $d.Paragraphs(5)
$insertAt = $d.Range($saifrPara.Range.Start, $saifrPara.Range.Start)
$insertAt.InsertBefore($newPara + "`r`r")

# ---------------------------------------------------------------------------
# 4) Small in-place text fixes inside the block that will be relocated below.
# ---------------------------------------------------------------------------

# 4a. "of  algebraic" (double space) -> "of algebraic"
$d.Content.Find.Execute(
    "framework of  algebraic quantum physics", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "framework of algebraic quantum physics", 2) | Out-Null

# 4b. Append the Radon-Nikodym sentence after "...quantum anomalies and information theory."
$d.Content.Find.Execute(
    "and the relationship this phenomena has with quantum anomalies and information theory.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and the relationship this phenomena has with quantum anomalies and information theory. In particular, the relationship between Radon-Nikodym cocycles and derivatives may provide a better understanding of the anomalous behavior of path integral measures through Tomita-Takesaki theory.",
    2) | Out-Null

# 4c. Append the WFIRST/EUCLID sentence after "...Precision Projector Laboratory."
$d.Content.Find.Execute(
    "In here I collaborated with JPL NASA researchers in the Precision Projector Laboratory.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In here I collaborated with JPL NASA researchers in the Precision Projector Laboratory. We studied motions in the centroid of images taken in simulated environments to improve the testing of infrared sensors. The sensors we studied are scheduled to be part of missions such as WFIRST and EUCLID.",
    2) | Out-Null

# 4d. Shorten the financial-request closing sentence.
$d.Content.Find.Execute(
    "purposes in Colombia has become increasingly difficult due to the low amount of the GDP invested in science.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "purposes in Colombia has become increasingly difficult.",
    2) | Out-Null

Write-Host "Done steps 1-4"
